$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for even_MAG-GUT27169.fa (row 6).
# All subsequent rows shift up by one.
$ws.Rows("6:6").Delete()

# After the above shift, the row for even_MAG-GUT38735.fa (originally row 9)
# is now at row 8. Delete it too; remaining rows shift up again.
$ws.Rows("8:8").Delete()
